$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.401.01"
$ws.Range("E2").Value = "  -1.35%  "

$ws.Range("D3").Value = "3.064.19"
$ws.Range("E3").Value = "  -2.72%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.06"
$ws.Range("E5").Value = "  -0.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.81"
$ws.Range("E6").Value = "  +4.61%  "

$ws.Range("E8").Value = "  +0.84%  "

$ws.Range("D9").Value = "3.061.69"
$ws.Range("E9").Value = "  -2.50%  "

$ws.Range("E10").Value = "  -4.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.82"
$ws.Range("E11").Value = "  -1.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.451"
$ws.Range("E12").Value = "  -1.73%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.94"
$ws.Range("E13").Value = "  -1.67%  "

$ws.Range("E14").Value = "  -4.24%  "

$ws.Range("E15").Value = "  -2.05%  "

$ws.Range("D16").Value = "3.570.04"
$ws.Range("E16").Value = "  -2.64%  "

$ws.Range("D17").Value = "63.426.01"
$ws.Range("E17").Value = "  -0.95%  "

$ws.Range("E18").Value = "  -2.10%  "

$ws.Range("D19").Value = "3.062.75"
$ws.Range("E19").Value = "  -2.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "471.69"
$ws.Range("E20").Value = "  +0.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.33"
$ws.Range("E21").Value = "  -0.87%  "

$ws.Range("E22").Value = "  -3.90%  "

$ws.Range("E23").Value = "  -1.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").Value = "  +1.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.62"
$ws.Range("E25").Value = "  -0.64%  "

$ws.Range("E26").Value = "  -2.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.37"
$ws.Range("E27").Value = "  +3.18%  "

$ws.Range("E28").Value = "  -0.28%  "

$ws.Range("E29").Value = "  +2.44%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("E31").Value = "  -2.32%  "

$ws.Range("E32").Value = "  -3.47%  "

$ws.Range("E33").Value = "  -2.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.12"
$ws.Range("E34").Value = "  -2.33%  "

$ws.Range("D35").Value = "0.0₃0820"
$ws.Range("E35").Value = "  -4.90%  "

$ws.Range("E36").Value = "  -1.99%  "

$ws.Range("E37").Value = "  +0.19%  "

$ws.Range("E38").Value = "  -3.28%  "

$ws.Range("E39").Value = "  -4.30%  "

$ws.Range("E40").Value = "  -1.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.20"
$ws.Range("E41").Value = "  -1.36%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "437.06"
$ws.Range("E42").Value = "  -5.96%  "

$ws.Range("E43").Value = "  -1.83%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.28"
$ws.Range("E44").Value = "  +3.62%  "

$ws.Range("E45").Value = "  +2.30%  "

$ws.Range("E46").Value = "  -4.46%  "

$ws.Range("D47").Value = "2.791.28"
$ws.Range("E47").Value = "  -3.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.71"
$ws.Range("E48").Value = "  -2.11%  "

$ws.Range("E49").Value = "  +0.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.98"
$ws.Range("E50").Value = "  +3.38%  "

$ws.Range("E51").Value = "  -0.38%  "
